# Auto-generated edit script applying the cell-value changes described in the diff.
# Each row updates currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ columns (H-N) per a scheduled price refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 4013.9614  # H17: 3920.8518 -> 4013.9614
$ws.Cells.Item(17, 10).Value = 4152  # J17: 4050 -> 4152
$ws.Cells.Item(17, 12).Value = 12456  # L17: 12150 -> 12456
$ws.Cells.Item(17, 14).Value = -12792  # N17: -12486 -> -12792

$ws.Cells.Item(32, 8).Value = 670.4  # H32: 604.7778 -> 670.4
$ws.Cells.Item(32, 9).Value = 550  # I32: 450.5 -> 550
$ws.Cells.Item(32, 10).Value = 750.6667  # J32: 648.8570999999999 -> 750.6667
$ws.Cells.Item(32, 11).Value = 550  # K32: 450.5 -> 550
$ws.Cells.Item(32, 12).Value = 750.6667  # L32: 648.8570999999999 -> 750.6667
$ws.Cells.Item(32, 13).Value = -224  # M32: -124.5 -> -224
$ws.Cells.Item(32, 14).Value = -1402.6667  # N32: -1300.8571 -> -1402.6667

$ws.Cells.Item(121, 8).Value = 900.8333  # H121: 898.8889 -> 900.8333
$ws.Cells.Item(121, 10).Value = 1176.25  # J121: 1055.7142 -> 1176.25
$ws.Cells.Item(121, 12).Value = 3528.75  # L121: 3167.1426 -> 3528.75
$ws.Cells.Item(121, 14).Value = -7022.75  # N121: -6661.142599999999 -> -7022.75

$ws.Cells.Item(132, 8).Value = 40875.848  # H132: 46049.78 -> 40875.848
$ws.Cells.Item(132, 9).Value = 42390.88  # I132: 46049.78 -> 42390.88
$ws.Cells.Item(132, 10).Value = 3000  # J132: 0 -> 3000
$ws.Cells.Item(132, 11).Value = 127172.64  # K132: 138149.34 -> 127172.64
$ws.Cells.Item(132, 12).Value = 9000  # L132: 0 -> 9000
$ws.Cells.Item(132, 13).Value = -124642.64  # M132: -135619.34 -> -124642.64
$ws.Cells.Item(132, 14).Value = -14060  # N132: None -> -14060

$ws.Cells.Item(137, 8).Value = 37038696  # H137: 19231796 -> 37038696
$ws.Cells.Item(137, 9).Value = 50001196  # I137: 23810256 -> 50001196
$ws.Cells.Item(137, 10).Value = 2983.2856  # J137: 2268.3 -> 2983.2856
$ws.Cells.Item(137, 11).Value = 150003588  # K137: 71430768 -> 150003588
$ws.Cells.Item(137, 12).Value = 8949.856800000001  # L137: 6804.900000000001 -> 8949.856800000001
$ws.Cells.Item(137, 13).Value = -150001038  # M137: -71428218 -> -150001038
$ws.Cells.Item(137, 14).Value = -14049.8568  # N137: -11904.9 -> -14049.8568

$ws.Cells.Item(138, 8).Value = 14709291  # H138: 14331977 -> 14709291
$ws.Cells.Item(138, 9).Value = 2804099.2  # I138: 2264951.5 -> 2804099.2
$ws.Cells.Item(138, 10).Value = 29415704  # J138: 38466028 -> 29415704
$ws.Cells.Item(138, 11).Value = 8412297.600000001  # K138: 6794854.5 -> 8412297.600000001
$ws.Cells.Item(138, 12).Value = 88247112  # L138: 115398084 -> 88247112
$ws.Cells.Item(138, 13).Value = -8407157.600000001  # M138: -6789714.5 -> -8407157.600000001
$ws.Cells.Item(138, 14).Value = -88257392  # N138: -115408364 -> -88257392

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(19, 8).Value = 536  # H19: 5000 -> 536
$ws.Cells.Item(19, 9).Value = 536  # I19: 0 -> 536
$ws.Cells.Item(19, 10).Value = 0  # J19: 5000 -> 0
$ws.Cells.Item(19, 11).Value = 536  # K19: 0 -> 536
$ws.Cells.Item(19, 12).Value = 0  # L19: 5000 -> 0
$ws.Cells.Item(19, 13).Value = -307  # M19: None -> -307
$ws.Cells.Item(19, 14).ClearContents()  # N19: -5458 -> (removed)

$ws.Cells.Item(74, 8).Value = 5116.788  # H74: 6307.077 -> 5116.788
$ws.Cells.Item(74, 9).Value = 990.1053000000001  # I74: 1111.7273 -> 990.1053000000001
$ws.Cells.Item(74, 10).Value = 10717.286  # J74: 10117 -> 10717.286
$ws.Cells.Item(74, 11).Value = 990.1053000000001  # K74: 1111.7273 -> 990.1053000000001
$ws.Cells.Item(74, 12).Value = 10717.286  # L74: 10117 -> 10717.286
$ws.Cells.Item(74, 13).Value = -116.1053000000001  # M74: -237.7273 -> -116.1053000000001
$ws.Cells.Item(74, 14).Value = -12465.286  # N74: -11865 -> -12465.286

$ws.Cells.Item(77, 8).Value = 5116.788  # H77: 6307.077 -> 5116.788
$ws.Cells.Item(77, 9).Value = 990.1053000000001  # I77: 1111.7273 -> 990.1053000000001
$ws.Cells.Item(77, 10).Value = 10717.286  # J77: 10117 -> 10717.286
$ws.Cells.Item(77, 11).Value = 4950.5265  # K77: 5558.636500000001 -> 4950.5265
$ws.Cells.Item(77, 12).Value = 53586.43  # L77: 50585 -> 53586.43
$ws.Cells.Item(77, 13).Value = -582.5264999999999  # M77: -1190.636500000001 -> -582.5264999999999
$ws.Cells.Item(77, 14).Value = -62322.43  # N77: -59321 -> -62322.43

$ws.Cells.Item(122, 8).Value = 2577.457  # H122: 2661.1765 -> 2577.457
$ws.Cells.Item(122, 9).Value = 1980.44  # I122: 2074.1667 -> 1980.44
$ws.Cells.Item(122, 11).Value = 5941.32  # K122: 6222.500100000001 -> 5941.32
$ws.Cells.Item(122, 13).Value = -3491.32  # M122: -3772.500100000001 -> -3491.32

$ws.Cells.Item(132, 8).Value = 3942.9167  # H132: 4417.381 -> 3942.9167
$ws.Cells.Item(132, 9).Value = 3497.389  # I132: 4001.7144 -> 3497.389
$ws.Cells.Item(132, 10).Value = 5279.5  # J132: 5248.7144 -> 5279.5
$ws.Cells.Item(132, 11).Value = 10492.167  # K132: 12005.1432 -> 10492.167
$ws.Cells.Item(132, 12).Value = 15838.5  # L132: 15746.1432 -> 15838.5
$ws.Cells.Item(132, 13).Value = -7962.167000000001  # M132: -9475.143199999999 -> -7962.167000000001
$ws.Cells.Item(132, 14).Value = -20898.5  # N132: -20806.1432 -> -20898.5

$ws.Cells.Item(139, 8).Value = 52190.832  # H139: 52536.25 -> 52190.832
$ws.Cells.Item(139, 10).Value = 52190.832  # J139: 52536.25 -> 52190.832
$ws.Cells.Item(139, 12).Value = 52190.832  # L139: 52536.25 -> 52190.832
$ws.Cells.Item(139, 14).Value = -62470.832  # N139: -62816.25 -> -62470.832

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(59, 8).Value = 54450  # H59: 54000 -> 54450
$ws.Cells.Item(59, 10).Value = 54450  # J59: 54000 -> 54450
$ws.Cells.Item(59, 12).Value = 54450  # L59: 54000 -> 54450
$ws.Cells.Item(59, 14).Value = -56144  # N59: -55694 -> -56144

$ws.Cells.Item(86, 8).Value = 1156.125  # H86: 1000 -> 1156.125
$ws.Cells.Item(86, 9).Value = 733  # I86: 1000 -> 733
$ws.Cells.Item(86, 10).Value = 1410  # J86: 0 -> 1410
$ws.Cells.Item(86, 11).Value = 733  # K86: 1000 -> 733
$ws.Cells.Item(86, 12).Value = 1410  # L86: 0 -> 1410
$ws.Cells.Item(86, 13).Value = 390  # M86: 123 -> 390
$ws.Cells.Item(86, 14).Value = -3656  # N86: None -> -3656

$ws.Cells.Item(89, 8).Value = 1156.125  # H89: 1000 -> 1156.125
$ws.Cells.Item(89, 9).Value = 733  # I89: 1000 -> 733
$ws.Cells.Item(89, 10).Value = 1410  # J89: 0 -> 1410
$ws.Cells.Item(89, 11).Value = 3665  # K89: 5000 -> 3665
$ws.Cells.Item(89, 12).Value = 7050  # L89: 0 -> 7050
$ws.Cells.Item(89, 13).Value = 1951  # M89: 616 -> 1951
$ws.Cells.Item(89, 14).Value = -18282  # N89: None -> -18282

$ws.Cells.Item(94, 8).Value = 1300.6296  # H94: 1125.421 -> 1300.6296
$ws.Cells.Item(94, 9).Value = 1384.8948  # I94: 873.875 -> 1384.8948
$ws.Cells.Item(94, 10).Value = 1100.5  # J94: 2467 -> 1100.5
$ws.Cells.Item(94, 11).Value = 1384.8948  # K94: 873.875 -> 1384.8948
$ws.Cells.Item(94, 12).Value = 1100.5  # L94: 2467 -> 1100.5
$ws.Cells.Item(94, 13).Value = -933.8948  # M94: -422.875 -> -933.8948
$ws.Cells.Item(94, 14).Value = -2002.5  # N94: -3369 -> -2002.5

$ws.Cells.Item(134, 8).Value = 3619.5117  # H134: 3823.3684 -> 3619.5117
$ws.Cells.Item(134, 9).Value = 1896.4667  # I134: 2136.2693 -> 1896.4667
$ws.Cells.Item(134, 10).Value = 7595.769  # J134: 7478.75 -> 7595.769
$ws.Cells.Item(134, 11).Value = 5689.4001  # K134: 6408.8079 -> 5689.4001
$ws.Cells.Item(134, 12).Value = 22787.307  # L134: 22436.25 -> 22787.307
$ws.Cells.Item(134, 13).Value = -3154.4001  # M134: -3873.8079 -> -3154.4001
$ws.Cells.Item(134, 14).Value = -27857.307  # N134: -27506.25 -> -27857.307

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2077.9756  # H58: 2133.2222 -> 2077.9756
$ws.Cells.Item(58, 9).Value = 1132.862  # I58: 1213.5555 -> 1132.862
$ws.Cells.Item(58, 10).Value = 4362  # J58: 4892.222 -> 4362
$ws.Cells.Item(58, 11).Value = 1132.862  # K58: 1213.5555 -> 1132.862
$ws.Cells.Item(58, 12).Value = 4362  # L58: 4892.222 -> 4362
$ws.Cells.Item(58, 13).Value = -929.8620000000001  # M58: -1010.5555 -> -929.8620000000001
$ws.Cells.Item(58, 14).Value = -4768  # N58: -5298.222 -> -4768

$ws.Cells.Item(136, 8).Value = 2077.9756  # H136: 2133.2222 -> 2077.9756
$ws.Cells.Item(136, 9).Value = 1132.862  # I136: 1213.5555 -> 1132.862
$ws.Cells.Item(136, 10).Value = 4362  # J136: 4892.222 -> 4362
$ws.Cells.Item(136, 11).Value = 3398.586  # K136: 3640.6665 -> 3398.586
$ws.Cells.Item(136, 12).Value = 13086  # L136: 14676.666 -> 13086
$ws.Cells.Item(136, 13).Value = -848.5860000000002  # M136: -1090.6665 -> -848.5860000000002
$ws.Cells.Item(136, 14).Value = -18186  # N136: -19776.666 -> -18186

$ws.Cells.Item(137, 8).Value = 36666.668  # H137: 42000 -> 36666.668
$ws.Cells.Item(137, 10).Value = 36666.668  # J137: 42000 -> 36666.668
$ws.Cells.Item(137, 12).Value = 36666.668  # L137: 42000 -> 36666.668
$ws.Cells.Item(137, 14).Value = -46866.668  # N137: -52200 -> -46866.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 196.95653  # H2: 165.96552 -> 196.95653
$ws.Cells.Item(2, 9).Value = 142.21428  # I2: 141.23077 -> 142.21428
$ws.Cells.Item(2, 10).Value = 282.1111  # J2: 186.0625 -> 282.1111
$ws.Cells.Item(2, 11).Value = 853.28568  # K2: 847.38462 -> 853.28568
$ws.Cells.Item(2, 12).Value = 1692.6666  # L2: 1116.375 -> 1692.6666
$ws.Cells.Item(2, 13).Value = -740.28568  # M2: -734.38462 -> -740.28568
$ws.Cells.Item(2, 14).Value = -1918.6666  # N2: -1342.375 -> -1918.6666

$ws.Cells.Item(19, 8).Value = 0  # H19: 902 -> 0
$ws.Cells.Item(19, 10).Value = 0  # J19: 902 -> 0
$ws.Cells.Item(19, 12).Value = 0  # L19: 2706 -> 0
$ws.Cells.Item(19, 14).ClearContents()  # N19: -3054 -> (removed)

$ws.Cells.Item(68, 8).Value = 12616.2  # H68: 11625 -> 12616.2
$ws.Cells.Item(68, 9).Value = 15636.5  # I68: 30625.5 -> 15636.5
$ws.Cells.Item(68, 10).Value = 535  # J68: 767.5714 -> 535
$ws.Cells.Item(68, 11).Value = 46909.5  # K68: 91876.5 -> 46909.5
$ws.Cells.Item(68, 12).Value = 1605  # L68: 2302.7142 -> 1605
$ws.Cells.Item(68, 13).Value = -46098.5  # M68: -91065.5 -> -46098.5
$ws.Cells.Item(68, 14).Value = -3227  # N68: -3924.7142 -> -3227

$ws.Cells.Item(71, 8).Value = 12616.2  # H71: 11625 -> 12616.2
$ws.Cells.Item(71, 9).Value = 15636.5  # I71: 30625.5 -> 15636.5
$ws.Cells.Item(71, 10).Value = 535  # J71: 767.5714 -> 535
$ws.Cells.Item(71, 11).Value = 140728.5  # K71: 275629.5 -> 140728.5
$ws.Cells.Item(71, 12).Value = 4815  # L71: 6908.1426 -> 4815
$ws.Cells.Item(71, 13).Value = -136672.5  # M71: -271573.5 -> -136672.5
$ws.Cells.Item(71, 14).Value = -12927  # N71: -15020.1426 -> -12927

$ws.Cells.Item(92, 8).Value = 2126.5715  # H92: 2142.889 -> 2126.5715
$ws.Cells.Item(92, 10).Value = 2221.5  # J92: 2183.7144 -> 2221.5
$ws.Cells.Item(92, 12).Value = 6664.5  # L92: 6551.1432 -> 6664.5
$ws.Cells.Item(92, 14).Value = -9160.5  # N92: -9047.143199999999 -> -9160.5

$ws.Cells.Item(107, 8).Value = 388.24445  # H107: 421.45715 -> 388.24445
$ws.Cells.Item(107, 9).Value = 400.53845  # I107: 418.09525 -> 400.53845
$ws.Cells.Item(107, 10).Value = 371.42105  # J107: 426.5 -> 371.42105
$ws.Cells.Item(107, 11).Value = 1201.61535  # K107: 1254.28575 -> 1201.61535
$ws.Cells.Item(107, 12).Value = 1114.26315  # L107: 1279.5 -> 1114.26315
$ws.Cells.Item(107, 13).Value = 718.38465  # M107: 665.71425 -> 718.38465
$ws.Cells.Item(107, 14).Value = -4954.26315  # N107: -5119.5 -> -4954.26315

$ws.Cells.Item(131, 8).Value = 8548684  # H131: 8773644 -> 8548684
$ws.Cells.Item(131, 10).Value = 9010746  # J131: 9261038 -> 9010746
$ws.Cells.Item(131, 12).Value = 27032238  # L131: 27783114 -> 27032238
$ws.Cells.Item(131, 14).Value = -27042318  # N131: -27793194 -> -27042318

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(6, 8).Value = 0  # H6: 29909 -> 0
$ws.Cells.Item(6, 10).Value = 0  # J6: 29909 -> 0
$ws.Cells.Item(6, 12).Value = 0  # L6: 29909 -> 0
$ws.Cells.Item(6, 14).ClearContents()  # N6: -30135 -> (removed)

$ws.Cells.Item(16, 8).Value = 0  # H16: 29909 -> 0
$ws.Cells.Item(16, 10).Value = 0  # J16: 29909 -> 0
$ws.Cells.Item(16, 12).Value = 0  # L16: 29909 -> 0
$ws.Cells.Item(16, 14).ClearContents()  # N16: -30409 -> (removed)

$ws.Cells.Item(70, 8).Value = 4915.625  # H70: 5049.8203 -> 4915.625
$ws.Cells.Item(70, 9).Value = 5104.5625  # I70: 5374.276 -> 5104.5625
$ws.Cells.Item(70, 10).Value = 4159.875  # J70: 4108.9 -> 4159.875
$ws.Cells.Item(70, 11).Value = 5104.5625  # K70: 5374.276 -> 5104.5625
$ws.Cells.Item(70, 12).Value = 4159.875  # L70: 4108.9 -> 4159.875
$ws.Cells.Item(70, 13).Value = -4834.5625  # M70: -5104.276 -> -4834.5625
$ws.Cells.Item(70, 14).Value = -4699.875  # N70: -4648.9 -> -4699.875

$ws.Cells.Item(73, 8).Value = 4915.625  # H73: 5049.8203 -> 4915.625
$ws.Cells.Item(73, 9).Value = 5104.5625  # I73: 5374.276 -> 5104.5625
$ws.Cells.Item(73, 10).Value = 4159.875  # J73: 4108.9 -> 4159.875
$ws.Cells.Item(73, 11).Value = 5104.5625  # K73: 5374.276 -> 5104.5625
$ws.Cells.Item(73, 12).Value = 4159.875  # L73: 4108.9 -> 4159.875
$ws.Cells.Item(73, 13).Value = -4168.5625  # M73: -4438.276 -> -4168.5625
$ws.Cells.Item(73, 14).Value = -6031.875  # N73: -5980.9 -> -6031.875

$ws.Cells.Item(102, 8).Value = 2148.7273  # H102: 3270.5 -> 2148.7273
$ws.Cells.Item(102, 9).Value = 2083.4707  # I102: 4240.125 -> 2083.4707
$ws.Cells.Item(102, 10).Value = 2218.0625  # J102: 2624.0833 -> 2218.0625
$ws.Cells.Item(102, 11).Value = 2083.4707  # K102: 4240.125 -> 2083.4707
$ws.Cells.Item(102, 12).Value = 2218.0625  # L102: 2624.0833 -> 2218.0625
$ws.Cells.Item(102, 13).Value = -461.4706999999999  # M102: -2618.125 -> -461.4706999999999
$ws.Cells.Item(102, 14).Value = -5462.0625  # N102: -5868.0833 -> -5462.0625

$ws.Cells.Item(126, 8).Value = 2699.7083  # H126: 2079.0527 -> 2699.7083
$ws.Cells.Item(126, 9).Value = 2149  # I126: 1481.9546 -> 2149
$ws.Cells.Item(126, 10).Value = 2975.0625  # J126: 2900.0625 -> 2975.0625
$ws.Cells.Item(126, 11).Value = 6447  # K126: 4445.8638 -> 6447
$ws.Cells.Item(126, 12).Value = 8925.1875  # L126: 8700.1875 -> 8925.1875
$ws.Cells.Item(126, 13).Value = -3977  # M126: -1975.8638 -> -3977
$ws.Cells.Item(126, 14).Value = -13865.1875  # N126: -13640.1875 -> -13865.1875

$ws.Cells.Item(132, 8).Value = 2268.4792  # H132: 2770.4055 -> 2268.4792
$ws.Cells.Item(132, 9).Value = 1857.125  # I132: 2203.6562 -> 1857.125
$ws.Cells.Item(132, 10).Value = 4325.25  # J132: 6397.6 -> 4325.25
$ws.Cells.Item(132, 11).Value = 5571.375  # K132: 6610.9686 -> 5571.375
$ws.Cells.Item(132, 12).Value = 12975.75  # L132: 19192.8 -> 12975.75
$ws.Cells.Item(132, 13).Value = -3041.375  # M132: -4080.9686 -> -3041.375
$ws.Cells.Item(132, 14).Value = -18035.75  # N132: -24252.8 -> -18035.75

$ws.Cells.Item(137, 8).Value = 45226.668  # H137: 46893.332 -> 45226.668
$ws.Cells.Item(137, 10).Value = 45226.668  # J137: 46893.332 -> 45226.668
$ws.Cells.Item(137, 12).Value = 45226.668  # L137: 46893.332 -> 45226.668
$ws.Cells.Item(137, 14).Value = -55426.668  # N137: -57093.332 -> -55426.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 2028.5714  # H46: 2116.6667 -> 2028.5714

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(7, 8).Value = 18000  # H7: 2375 -> 18000
$ws.Cells.Item(7, 9).Value = 25500  # I7: 500 -> 25500
$ws.Cells.Item(7, 11).Value = 25500  # K7: 500 -> 25500
$ws.Cells.Item(7, 13).Value = -25387  # M7: -387 -> -25387

$ws.Cells.Item(132, 8).Value = 2301.2563  # H132: 2210.0789 -> 2301.2563
$ws.Cells.Item(132, 9).Value = 2077.8215  # I132: 2033.8889 -> 2077.8215
$ws.Cells.Item(132, 10).Value = 2870  # J132: 2642.5454 -> 2870
$ws.Cells.Item(132, 11).Value = 6233.4645  # K132: 6101.6667 -> 6233.4645
$ws.Cells.Item(132, 12).Value = 8610  # L132: 7927.6362 -> 8610
$ws.Cells.Item(132, 13).Value = -3703.4645  # M132: -3571.6667 -> -3703.4645
$ws.Cells.Item(132, 14).Value = -13670  # N132: -12987.6362 -> -13670

$ws.Cells.Item(136, 8).Value = 1678.0385  # H136: 1861.1305 -> 1678.0385
$ws.Cells.Item(136, 9).Value = 992.2105  # I136: 1029.3889 -> 992.2105
$ws.Cells.Item(136, 10).Value = 3539.5715  # J136: 4855.4 -> 3539.5715
$ws.Cells.Item(136, 11).Value = 2976.6315  # K136: 3088.1667 -> 2976.6315
$ws.Cells.Item(136, 12).Value = 10618.7145  # L136: 14566.2 -> 10618.7145
$ws.Cells.Item(136, 13).Value = -426.6315  # M136: -538.1666999999998 -> -426.6315
$ws.Cells.Item(136, 14).Value = -15718.7145  # N136: -19666.2 -> -15718.7145
